$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) from 45182 to 45184 for the data rows (2..271)
$lastRow = 271
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45184
